$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.040629009448294
$ws.Range("D2").Value = 1.048480843738454
$ws.Range("E2").Value = 1.038992835535349
$ws.Range("F2").Value = 1.057299806218056
$ws.Range("I2").Value = 1.041110595859552
$ws.Range("J2").Value = 1.045714584523252
$ws.Range("K2").Value = 1.051240569656995
$ws.Range("L2").Value = 1.041779289609118
$ws.Range("M2").Value = 1.060035154805782
$ws.Range("N2").Value = 1.019094698327555
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.041645512437985
$ws.Range("D3").Value = 1.049319489136312
$ws.Range("E3").Value = 1.039859850047894
$ws.Range("F3").Value = 1.058371419347555
$ws.Range("I3").Value = 1.041389751084952
$ws.Range("J3").Value = 1.046376449235785
$ws.Range("K3").Value = 1.051891200176147
$ws.Range("L3").Value = 1.042456281235588
$ws.Range("M3").Value = 1.060919921848696
$ws.Range("N3").Value = 1.019319259435051
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.042303519128701
$ws.Range("D4").Value = 1.049862378553369
$ws.Range("E4").Value = 1.040421429472954
$ws.Range("F4").Value = 1.059065559865469
$ws.Range("I4").Value = 1.04156932495967
$ws.Range("J4").Value = 1.046804386636084
$ws.Range("K4").Value = 1.052311784894145
$ws.Range("L4").Value = 1.04289426401213
$ws.Range("M4").Value = 1.061492560257697
$ws.Range("N4").Value = 1.019464337280788
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.042580207416925
$ws.Range("D5").Value = 1.050090663497465
$ws.Range("E5").Value = 1.040657651446501
$ws.Range("F5").Value = 1.059357551955204
$ws.Range("I5").Value = 1.041644564178026
$ws.Range("J5").Value = 1.046984211223203
$ws.Range("K5").Value = 1.052488498404961
$ws.Range("L5").Value = 1.043078373170071
$ws.Range("M5").Value = 1.061733329199949
$ws.Range("N5").Value = 1.019525273109648
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.042626668241948
$ws.Range("D6").Value = 1.050128996701901
$ws.Range("E6").Value = 1.040697321988487
$ws.Range("F6").Value = 1.059406588981489
$ws.Range("I6").Value = 1.041657182291761
$ws.Range("J6").Value = 1.047014399839662
$ws.Range("K6").Value = 1.052518163483519
$ws.Range("L6").Value = 1.043109284788214
$ws.Range("M6").Value = 1.061773757210817
$ws.Range("N6").Value = 1.019535501276552
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.04230721600678
$ws.Range("D7").Value = 1.049865428695828
$ws.Range("E7").Value = 1.040424585357079
$ws.Range("F7").Value = 1.059069460789517
$ws.Range("I7").Value = 1.04157033130678
$ws.Range("J7").Value = 1.046806789777141
$ws.Range("K7").Value = 1.05231414654352
$ws.Range("L7").Value = 1.042896724162317
$ws.Range("M7").Value = 1.061495777299999
$ws.Range("N7").Value = 1.019465151724355
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.040972487124638
$ws.Range("D8").Value = 1.048764219837441
$ws.Range("E8").Value = 1.039285730074144
$ws.Range("F8").Value = 1.05766181031913
$ws.Range("I8").Value = 1.041205156607286
$ws.Range("J8").Value = 1.045938333519858
$ws.Range("K8").Value = 1.05146053925109
$ws.Range("L8").Value = 1.042008097496459
$ws.Range("M8").Value = 1.06033413768848
$ws.Range("N8").Value = 1.019170636898053
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.03862253741455
$ws.Range("D9").Value = 1.046825542831023
$ws.Range("E9").Value = 1.037283269711651
$ws.Range("F9").Value = 1.055187006447214
$ws.Range("I9").Value = 1.040553580514101
$ws.Range("J9").Value = 1.044405471535625
$ws.Range("K9").Value = 1.049953206039619
$ws.Range("L9").Value = 1.040441666973413
$ws.Range("M9").Value = 1.058288238395053
$ws.Range("N9").Value = 1.018649927954435
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.037057271515743
$ws.Range("D10").Value = 1.045534341719953
$ws.Range("E10").Value = 1.035951266359435
$ws.Range("F10").Value = 1.05354097246383
$ws.Range("I10").Value = 1.040113772980993
$ws.Range("J10").Value = 1.043381889655195
$ws.Range("K10").Value = 1.048946217541588
$ws.Range("L10").Value = 1.039397039062416
$ws.Range("M10").Value = 1.05692504802395
$ws.Range("N10").Value = 1.018301636494413
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.036379820490628
$ws.Range("D11").Value = 1.044975542744974
$ws.Range("E11").Value = 1.035375208083623
$ws.Range("F11").Value = 1.052829137395936
$ws.Range("I11").Value = 1.039922048370014
$ws.Range("J11").Value = 1.042938275760044
$ws.Range("K11").Value = 1.048509688309415
$ws.Range("L11").Value = 1.038944629393029
$ws.Range("M11").Value = 1.05633495359602
$ws.Range("N11").Value = 1.018150552202295
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.036128233077782
$ws.Range("D12").Value = 1.044768025589854
$ws.Range("E12").Value = 1.035161341456223
$ws.Range("F12").Value = 1.052564867064516
$ws.Range("I12").Value = 1.039850640465739
$ws.Range("J12").Value = 1.04277343865214
$ws.Range("K12").Value = 1.048347467723978
$ws.Range("L12").Value = 1.038776572919779
$ws.Range("M12").Value = 1.056115793095341
$ws.Range("N12").Value = 1.018094392157339
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.036182197255828
$ws.Range("D13").Value = 1.044812536647203
$ws.Range("E13").Value = 1.035207211717593
$ws.Range("F13").Value = 1.052621547745666
$ws.Range("I13").Value = 1.039865966426846
$ws.Range("J13").Value = 1.042808799442325
$ws.Range("K13").Value = 1.048382267940947
$ws.Range("L13").Value = 1.038812622111051
$ws.Range("M13").Value = 1.056162802545468
$ws.Range("N13").Value = 1.018106440511832
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.03635902321574
$ws.Range("D14").Value = 1.044958388376891
$ws.Range("E14").Value = 1.035357527613328
$ws.Range("F14").Value = 1.052807289932626
$ws.Range("I14").Value = 1.039916149703134
$ws.Range("J14").Value = 1.042924651477902
$ws.Range("K14").Value = 1.048496280616326
$ws.Range("L14").Value = 1.038930738021896
$ws.Range("M14").Value = 1.056316837170723
$ws.Range("N14").Value = 1.018145910822645
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.036467977896408
$ws.Range("D15").Value = 1.045048258511728
$ws.Range("E15").Value = 1.035450156405922
$ws.Range("F15").Value = 1.052921749978274
$ws.Range("I15").Value = 1.039947043728309
$ws.Range("J15").Value = 1.042996023912992
$ws.Range("K15").Value = 1.048566517769695
$ws.Range("L15").Value = 1.039003511639328
$ws.Range("M15").Value = 1.056411746556281
$ws.Range("N15").Value = 1.018170224408997
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.03710223844515
$ws.Range("D16").Value = 1.045571433735736
$ws.Range("E16").Value = 1.035989512425872
$ws.Range("F16").Value = 1.053588233789408
$ws.Range("I16").Value = 1.040126470055957
$ws.Range("J16").Value = 1.043411322532507
$ws.Range("K16").Value = 1.048975178122019
$ws.Range("L16").Value = 1.039427062380483
$ws.Range("M16").Value = 1.056964214390027
$ws.Range("N16").Value = 1.018311657753889
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.037500178900662
$ws.Range("D17").Value = 1.045899688408211
$ws.Range("E17").Value = 1.036328026245062
$ws.Range("F17").Value = 1.054006545000757
$ws.Range("I17").Value = 1.040238675559473
$ws.Range("J17").Value = 1.043671722572737
$ws.Range("K17").Value = 1.049231387309028
$ws.Range("L17").Value = 1.039692723783298
$ws.Range("M17").Value = 1.057310810268782
$ws.Range("N17").Value = 1.018400302544759
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.037732321651992
$ws.Range("D18").Value = 1.046091182671083
$ws.Range("E18").Value = 1.036525543889943
$ws.Range("F18").Value = 1.05425062648307
$ws.Range("I18").Value = 1.040303999128258
$ws.Range("J18").Value = 1.04382357125961
$ws.Range("K18").Value = 1.049380781853826
$ws.Range("L18").Value = 1.039847671941291
$ws.Range("M18").Value = 1.057512990735089
$ws.Range("N18").Value = 1.018451981312707
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.037811481556018
$ws.Range("D19").Value = 1.046156482142514
$ws.Range("E19").Value = 1.036592903859511
$ws.Range("F19").Value = 1.05433386681193
$ws.Range("I19").Value = 1.040326251737482
$ws.Range("J19").Value = 1.043875341238544
$ws.Range("K19").Value = 1.049431713389082
$ws.Range("L19").Value = 1.039900503948056
$ws.Range("M19").Value = 1.057581931912824
$ws.Range("N19").Value = 1.018469597979772
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.037457480467793
$ws.Range("D20").Value = 1.045864466797927
$ws.Range("E20").Value = 1.036291699850213
$ws.Range("F20").Value = 1.053961655099585
$ws.Range("I20").Value = 1.040226649795472
$ws.Range("J20").Value = 1.043643788066589
$ws.Range("K20").Value = 1.049203903443557
$ws.Range("L20").Value = 1.039664221616111
$ws.Range("M20").Value = 1.057273622054519
$ws.Range("N20").Value = 1.018390794510008
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.036306951050504
$ws.Range("D21").Value = 1.044915437411966
$ws.Range("E21").Value = 1.035313260367397
$ws.Range("F21").Value = 1.052752589705413
$ws.Range("I21").Value = 1.039901377300827
$ws.Range("J21").Value = 1.042890537570976
$ws.Range("K21").Value = 1.048462708769868
$ws.Range("L21").Value = 1.038895956140402
$ws.Range("M21").Value = 1.05627147706853
$ws.Range("N21").Value = 1.018134288919186
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.035583845648112
$ws.Range("D22").Value = 1.044319009494487
$ws.Range("E22").Value = 1.034698696424714
$ws.Range("F22").Value = 1.051993194615445
$ws.Range("I22").Value = 1.039695749928029
$ws.Range("J22").Value = 1.042416597016979
$ws.Range("K22").Value = 1.047996261213216
$ws.Range("L22").Value = 1.038412851546746
$ws.Range("M22").Value = 1.055641544252407
$ws.Range("N22").Value = 1.017972778815897
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.035967151117709
$ws.Range("D23").Value = 1.044635161927372
$ws.Range("E23").Value = 1.035024429414901
$ws.Range("F23").Value = 1.052395689079686
$ws.Range("I23").Value = 1.039804862566248
$ws.Range("J23").Value = 1.042667874173071
$ws.Range("K23").Value = 1.048243574419198
$ws.Range("L23").Value = 1.038668960491662
$ws.Range("M23").Value = 1.055975468681132
$ws.Range("N23").Value = 1.018058420565235
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.037476773956346
$ws.Range("D24").Value = 1.045880381838332
$ws.Range("E24").Value = 1.036308113973495
$ws.Range("F24").Value = 1.053981938642847
$ws.Range("I24").Value = 1.040232084103712
$ws.Range("J24").Value = 1.043656410585649
$ws.Range("K24").Value = 1.049216322367435
$ws.Range("L24").Value = 1.039677100542209
$ws.Range("M24").Value = 1.057290425755367
$ws.Range("N24").Value = 1.018395090862721
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.039229815778998
$ws.Range("D25").Value = 1.047326520188784
$ws.Range("E25").Value = 1.037800433618193
$ws.Range("F25").Value = 1.055826128962939
$ws.Range("I25").Value = 1.040722985720171
$ws.Range("J25").Value = 1.044405471535625
$ws.Range("K25").Value = 1.049953206039619
$ws.Range("L25").Value = 1.040441666973413
$ws.Range("M25").Value = 1.058288238395053
$ws.Range("N25").Value = 1.018649927954435
